# Auto-generated Excel COM-interop script
# Applies the brazil_serie-b_2023 row corrections described by the commit diff:
#  - 12 groups of adjacent match rows had their match data (home..url, columns F:V)
#    rotated among themselves (rows 64/65, 74/76, 85/86, 100/102, 112/113, 115/116,
#    143/144, 163/164, 166/167, 177/178/179, 203/204, 205/206).
#  - One brand-new match row (282) was appended at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# rotate rows (64, 65)
$ws.Cells.Item(64,6).Value2 = "Sampaio Correa"
$ws.Cells.Item(64,7).Value2 = 1
$ws.Cells.Item(64,8).Value2 = "ABC"
$ws.Cells.Item(64,9).Value2 = 0
$ws.Cells.Item(64,10).Value2 = 1.93
$ws.Cells.Item(64,11).Value2 = "14/05/2023 23:12"
$ws.Cells.Item(64,12).Value2 = 1.91
$ws.Cells.Item(64,13).Value2 = "20/05/2023 21:52"
$ws.Cells.Item(64,14).Value2 = 3.34
$ws.Cells.Item(64,15).Value2 = "14/05/2023 23:12"
$ws.Cells.Item(64,16).Value2 = 3.47
$ws.Cells.Item(64,17).Value2 = "20/05/2023 21:57"
$ws.Cells.Item(64,18).Value2 = 4.4
$ws.Cells.Item(64,19).Value2 = "14/05/2023 23:12"
$ws.Cells.Item(64,20).Value2 = 4.47
$ws.Cells.Item(64,21).Value2 = "20/05/2023 21:57"
$ws.Cells.Item(64,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/sampaio-correa-abc/f77MOPo4/"
$ws.Cells.Item(65,6).Value2 = "Ituano"
$ws.Cells.Item(65,7).Value2 = 0
$ws.Cells.Item(65,8).Value2 = "Novorizontino"
$ws.Cells.Item(65,9).Value2 = 2
$ws.Cells.Item(65,10).Value2 = 2.65
$ws.Cells.Item(65,11).Value2 = "15/05/2023 01:42"
$ws.Cells.Item(65,12).Value2 = 2.71
$ws.Cells.Item(65,13).Value2 = "20/05/2023 21:59"
$ws.Cells.Item(65,14).Value2 = 3.13
$ws.Cells.Item(65,15).Value2 = "15/05/2023 01:42"
$ws.Cells.Item(65,16).Value2 = 3.08
$ws.Cells.Item(65,17).Value2 = "20/05/2023 21:59"
$ws.Cells.Item(65,18).Value2 = 2.75
$ws.Cells.Item(65,19).Value2 = "15/05/2023 01:42"
$ws.Cells.Item(65,20).Value2 = 2.94
$ws.Cells.Item(65,21).Value2 = "20/05/2023 21:50"
$ws.Cells.Item(65,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/ituano-novorizontino/KlF9RR1o/"

# rotate rows (74, 76)
$ws.Cells.Item(74,6).Value2 = "Vila Nova FC"
$ws.Cells.Item(74,7).Value2 = 1
$ws.Cells.Item(74,8).Value2 = "Ituano"
$ws.Cells.Item(74,9).Value2 = 0
$ws.Cells.Item(74,10).Value2 = 1.84
$ws.Cells.Item(74,11).Value2 = "21/05/2023 23:42"
$ws.Cells.Item(74,12).Value2 = 1.77
$ws.Cells.Item(74,13).Value2 = "24/05/2023 23:56"
$ws.Cells.Item(74,14).Value2 = 3.32
$ws.Cells.Item(74,15).Value2 = "21/05/2023 23:42"
$ws.Cells.Item(74,16).Value2 = 3.39
$ws.Cells.Item(74,17).Value2 = "24/05/2023 23:53"
$ws.Cells.Item(74,18).Value2 = 4.98
$ws.Cells.Item(74,19).Value2 = "21/05/2023 23:42"
$ws.Cells.Item(74,20).Value2 = 5.65
$ws.Cells.Item(74,21).Value2 = "24/05/2023 23:56"
$ws.Cells.Item(74,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/vila-nova-fc-ituano/428cG1gp/"
$ws.Cells.Item(76,6).Value2 = "Londrina"
$ws.Cells.Item(76,7).Value2 = 1
$ws.Cells.Item(76,8).Value2 = "Ceara"
$ws.Cells.Item(76,9).Value2 = 3
$ws.Cells.Item(76,10).Value2 = 2.67
$ws.Cells.Item(76,11).Value2 = "21/05/2023 20:42"
$ws.Cells.Item(76,12).Value2 = 2.83
$ws.Cells.Item(76,13).Value2 = "24/05/2023 23:59"
$ws.Cells.Item(76,14).Value2 = 2.98
$ws.Cells.Item(76,15).Value2 = "21/05/2023 20:42"
$ws.Cells.Item(76,16).Value2 = 3.12
$ws.Cells.Item(76,17).Value2 = "24/05/2023 23:59"
$ws.Cells.Item(76,18).Value2 = 3
$ws.Cells.Item(76,19).Value2 = "21/05/2023 20:42"
$ws.Cells.Item(76,20).Value2 = 2.79
$ws.Cells.Item(76,21).Value2 = "24/05/2023 23:58"
$ws.Cells.Item(76,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/londrina-ceara/pAPTIfFr/"

# rotate rows (85, 86)
$ws.Cells.Item(85,6).Value2 = "Ceara"
$ws.Cells.Item(85,7).Value2 = 0
$ws.Cells.Item(85,8).Value2 = "Novorizontino"
$ws.Cells.Item(85,9).Value2 = 3
$ws.Cells.Item(85,10).Value2 = 2.05
$ws.Cells.Item(85,11).Value2 = "26/05/2023 03:42"
$ws.Cells.Item(85,12).Value2 = 2.09
$ws.Cells.Item(85,13).Value2 = "28/05/2023 20:20"
$ws.Cells.Item(85,14).Value2 = 3.22
$ws.Cells.Item(85,15).Value2 = "26/05/2023 03:42"
$ws.Cells.Item(85,16).Value2 = 3.18
$ws.Cells.Item(85,17).Value2 = "28/05/2023 20:29"
$ws.Cells.Item(85,18).Value2 = 3.8
$ws.Cells.Item(85,19).Value2 = "26/05/2023 03:42"
$ws.Cells.Item(85,20).Value2 = 4.15
$ws.Cells.Item(85,21).Value2 = "28/05/2023 20:29"
$ws.Cells.Item(85,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/ceara-novorizontino/2BcuyfyR/"
$ws.Cells.Item(86,6).Value2 = "Mirassol"
$ws.Cells.Item(86,7).Value2 = 1
$ws.Cells.Item(86,8).Value2 = "Criciuma"
$ws.Cells.Item(86,9).Value2 = 0
$ws.Cells.Item(86,10).Value2 = 2.04
$ws.Cells.Item(86,11).Value2 = "26/05/2023 03:42"
$ws.Cells.Item(86,12).Value2 = 2
$ws.Cells.Item(86,13).Value2 = "28/05/2023 20:21"
$ws.Cells.Item(86,14).Value2 = 3.36
$ws.Cells.Item(86,15).Value2 = "26/05/2023 03:42"
$ws.Cells.Item(86,16).Value2 = 3.31
$ws.Cells.Item(86,17).Value2 = "28/05/2023 20:21"
$ws.Cells.Item(86,18).Value2 = 4.05
$ws.Cells.Item(86,19).Value2 = "26/05/2023 03:42"
$ws.Cells.Item(86,20).Value2 = 4.3
$ws.Cells.Item(86,21).Value2 = "28/05/2023 20:21"
$ws.Cells.Item(86,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/mirassol-criciuma/td9WxG6E/"

# rotate rows (100, 102)
$ws.Cells.Item(100,6).Value2 = "Tombense"
$ws.Cells.Item(100,7).Value2 = 1
$ws.Cells.Item(100,8).Value2 = "Vitoria"
$ws.Cells.Item(100,9).Value2 = 2
$ws.Cells.Item(100,10).Value2 = 2.46
$ws.Cells.Item(100,11).Value2 = "03/06/2023 02:42"
$ws.Cells.Item(100,12).Value2 = 2.46
$ws.Cells.Item(100,13).Value2 = "06/06/2023 23:59"
$ws.Cells.Item(100,14).Value2 = 3.03
$ws.Cells.Item(100,15).Value2 = "03/06/2023 02:42"
$ws.Cells.Item(100,16).Value2 = 3.18
$ws.Cells.Item(100,17).Value2 = "06/06/2023 23:59"
$ws.Cells.Item(100,18).Value2 = 3.08
$ws.Cells.Item(100,19).Value2 = "03/06/2023 02:42"
$ws.Cells.Item(100,20).Value2 = 3.2
$ws.Cells.Item(100,21).Value2 = "06/06/2023 23:59"
$ws.Cells.Item(100,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/tombense-vitoria/Kp3SflSa/"
$ws.Cells.Item(102,6).Value2 = "Chapecoense-SC"
$ws.Cells.Item(102,7).Value2 = 0
$ws.Cells.Item(102,8).Value2 = "Vila Nova FC"
$ws.Cells.Item(102,9).Value2 = 1
$ws.Cells.Item(102,10).Value2 = 2.37
$ws.Cells.Item(102,11).Value2 = "04/06/2023 00:13"
$ws.Cells.Item(102,12).Value2 = 2.81
$ws.Cells.Item(102,13).Value2 = "06/06/2023 23:59"
$ws.Cells.Item(102,14).Value2 = 3.05
$ws.Cells.Item(102,15).Value2 = "04/06/2023 00:13"
$ws.Cells.Item(102,16).Value2 = 2.99
$ws.Cells.Item(102,17).Value2 = "06/06/2023 23:59"
$ws.Cells.Item(102,18).Value2 = 3.22
$ws.Cells.Item(102,19).Value2 = "04/06/2023 00:13"
$ws.Cells.Item(102,20).Value2 = 2.92
$ws.Cells.Item(102,21).Value2 = "06/06/2023 23:59"
$ws.Cells.Item(102,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/chapecoense-sc-vila-nova-fc/6saOeUCg/"

# rotate rows (112, 113)
$ws.Cells.Item(112,6).Value2 = "Novorizontino"
$ws.Cells.Item(112,7).Value2 = 1
$ws.Cells.Item(112,8).Value2 = "Sampaio Correa"
$ws.Cells.Item(112,9).Value2 = 0
$ws.Cells.Item(112,10).Value2 = 1.56
$ws.Cells.Item(112,11).Value2 = "08/06/2023 14:42"
$ws.Cells.Item(112,12).Value2 = 1.61
$ws.Cells.Item(112,13).Value2 = "10/06/2023 21:22"
$ws.Cells.Item(112,14).Value2 = 3.64
$ws.Cells.Item(112,15).Value2 = "08/06/2023 14:42"
$ws.Cells.Item(112,16).Value2 = 3.81
$ws.Cells.Item(112,17).Value2 = "10/06/2023 21:22"
$ws.Cells.Item(112,18).Value2 = 6.5
$ws.Cells.Item(112,19).Value2 = "08/06/2023 14:42"
$ws.Cells.Item(112,20).Value2 = 6.55
$ws.Cells.Item(112,21).Value2 = "10/06/2023 21:22"
$ws.Cells.Item(112,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/novorizontino-sampaio-correa/QNPgB9Ct/"
$ws.Cells.Item(113,6).Value2 = "Ituano"
$ws.Cells.Item(113,7).Value2 = 1
$ws.Cells.Item(113,8).Value2 = "Atletico GO"
$ws.Cells.Item(113,9).Value2 = 1
$ws.Cells.Item(113,10).Value2 = 2.51
$ws.Cells.Item(113,11).Value2 = "07/06/2023 02:42"
$ws.Cells.Item(113,12).Value2 = 2.28
$ws.Cells.Item(113,13).Value2 = "10/06/2023 21:51"
$ws.Cells.Item(113,14).Value2 = 3.03
$ws.Cells.Item(113,15).Value2 = "07/06/2023 02:42"
$ws.Cells.Item(113,16).Value2 = 3.22
$ws.Cells.Item(113,17).Value2 = "10/06/2023 21:59"
$ws.Cells.Item(113,18).Value2 = 3.02
$ws.Cells.Item(113,19).Value2 = "07/06/2023 02:42"
$ws.Cells.Item(113,20).Value2 = 3.5
$ws.Cells.Item(113,21).Value2 = "10/06/2023 21:59"
$ws.Cells.Item(113,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/ituano-atletico-go/jVYRFk4P/"

# rotate rows (115, 116)
$ws.Cells.Item(115,6).Value2 = "Londrina"
$ws.Cells.Item(115,7).Value2 = 1
$ws.Cells.Item(115,8).Value2 = "Mirassol"
$ws.Cells.Item(115,9).Value2 = 2
$ws.Cells.Item(115,10).Value2 = 2.93
$ws.Cells.Item(115,11).Value2 = "08/06/2023 14:42"
$ws.Cells.Item(115,12).Value2 = 3.75
$ws.Cells.Item(115,13).Value2 = "11/06/2023 15:59"
$ws.Cells.Item(115,14).Value2 = 3.03
$ws.Cells.Item(115,15).Value2 = "08/06/2023 14:42"
$ws.Cells.Item(115,16).Value2 = 3.12
$ws.Cells.Item(115,17).Value2 = "11/06/2023 15:59"
$ws.Cells.Item(115,18).Value2 = 2.57
$ws.Cells.Item(115,19).Value2 = "08/06/2023 14:42"
$ws.Cells.Item(115,20).Value2 = 2.24
$ws.Cells.Item(115,21).Value2 = "11/06/2023 15:59"
$ws.Cells.Item(115,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/londrina-mirassol/f3HA7RB5/"
$ws.Cells.Item(116,6).Value2 = "Ponte Preta"
$ws.Cells.Item(116,7).Value2 = 1
$ws.Cells.Item(116,8).Value2 = "Sport Recife"
$ws.Cells.Item(116,9).Value2 = 1
$ws.Cells.Item(116,10).Value2 = 3.04
$ws.Cells.Item(116,11).Value2 = "08/06/2023 14:42"
$ws.Cells.Item(116,12).Value2 = 3.72
$ws.Cells.Item(116,13).Value2 = "11/06/2023 15:53"
$ws.Cells.Item(116,14).Value2 = 3
$ws.Cells.Item(116,15).Value2 = "08/06/2023 14:42"
$ws.Cells.Item(116,16).Value2 = 3.11
$ws.Cells.Item(116,17).Value2 = "11/06/2023 15:53"
$ws.Cells.Item(116,18).Value2 = 2.52
$ws.Cells.Item(116,19).Value2 = "08/06/2023 14:42"
$ws.Cells.Item(116,20).Value2 = 2.26
$ws.Cells.Item(116,21).Value2 = "11/06/2023 15:53"
$ws.Cells.Item(116,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/ponte-preta-sport-recife/t0iygSdC/"

# rotate rows (143, 144)
$ws.Cells.Item(143,6).Value2 = "Sampaio Correa"
$ws.Cells.Item(143,7).Value2 = 0
$ws.Cells.Item(143,8).Value2 = "Criciuma"
$ws.Cells.Item(143,9).Value2 = 2
$ws.Cells.Item(143,10).Value2 = 2.48
$ws.Cells.Item(143,11).Value2 = "29/06/2023 02:42"
$ws.Cells.Item(143,12).Value2 = 2.16
$ws.Cells.Item(143,13).Value2 = "01/07/2023 21:46"
$ws.Cells.Item(143,14).Value2 = 2.97
$ws.Cells.Item(143,15).Value2 = "29/06/2023 02:42"
$ws.Cells.Item(143,16).Value2 = 3.04
$ws.Cells.Item(143,17).Value2 = "01/07/2023 21:46"
$ws.Cells.Item(143,18).Value2 = 3.12
$ws.Cells.Item(143,19).Value2 = "29/06/2023 02:42"
$ws.Cells.Item(143,20).Value2 = 4.15
$ws.Cells.Item(143,21).Value2 = "01/07/2023 21:46"
$ws.Cells.Item(143,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/sampaio-correa-criciuma/dUS2YfGL/"
$ws.Cells.Item(144,6).Value2 = "Avai"
$ws.Cells.Item(144,7).Value2 = 0
$ws.Cells.Item(144,8).Value2 = "ABC"
$ws.Cells.Item(144,9).Value2 = 2
$ws.Cells.Item(144,10).Value2 = 1.99
$ws.Cells.Item(144,11).Value2 = "29/06/2023 00:12"
$ws.Cells.Item(144,12).Value2 = 1.99
$ws.Cells.Item(144,13).Value2 = "01/07/2023 21:58"
$ws.Cells.Item(144,14).Value2 = 3.19
$ws.Cells.Item(144,15).Value2 = "29/06/2023 00:12"
$ws.Cells.Item(144,16).Value2 = 3.11
$ws.Cells.Item(144,17).Value2 = "01/07/2023 21:58"
$ws.Cells.Item(144,18).Value2 = 4.39
$ws.Cells.Item(144,19).Value2 = "29/06/2023 00:12"
$ws.Cells.Item(144,20).Value2 = 4.75
$ws.Cells.Item(144,21).Value2 = "01/07/2023 21:58"
$ws.Cells.Item(144,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/avai-abc/88n8CDg2/"

# rotate rows (163, 164)
$ws.Cells.Item(163,6).Value2 = "Sampaio Correa"
$ws.Cells.Item(163,7).Value2 = 0
$ws.Cells.Item(163,8).Value2 = "Ituano"
$ws.Cells.Item(163,9).Value2 = 0
$ws.Cells.Item(163,10).Value2 = 1.99
$ws.Cells.Item(163,11).Value2 = "08/07/2023 23:12"
$ws.Cells.Item(163,12).Value2 = 1.95
$ws.Cells.Item(163,13).Value2 = "15/07/2023 02:21"
$ws.Cells.Item(163,14).Value2 = 3.21
$ws.Cells.Item(163,15).Value2 = "08/07/2023 23:12"
$ws.Cells.Item(163,16).Value2 = 3.2
$ws.Cells.Item(163,17).Value2 = "15/07/2023 02:20"
$ws.Cells.Item(163,18).Value2 = 4.35
$ws.Cells.Item(163,19).Value2 = "08/07/2023 23:12"
$ws.Cells.Item(163,20).Value2 = 4.74
$ws.Cells.Item(163,21).Value2 = "15/07/2023 02:21"
$ws.Cells.Item(163,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/sampaio-correa-ituano/KnWlRXPh/"
$ws.Cells.Item(164,6).Value2 = "Atletico GO"
$ws.Cells.Item(164,7).Value2 = 3
$ws.Cells.Item(164,8).Value2 = "Sport Recife"
$ws.Cells.Item(164,9).Value2 = 1
$ws.Cells.Item(164,10).Value2 = 2.64
$ws.Cells.Item(164,11).Value2 = "09/07/2023 23:11"
$ws.Cells.Item(164,12).Value2 = 2.68
$ws.Cells.Item(164,13).Value2 = "15/07/2023 02:27"
$ws.Cells.Item(164,14).Value2 = 3.01
$ws.Cells.Item(164,15).Value2 = "09/07/2023 23:11"
$ws.Cells.Item(164,16).Value2 = 2.99
$ws.Cells.Item(164,17).Value2 = "15/07/2023 02:27"
$ws.Cells.Item(164,18).Value2 = 2.87
$ws.Cells.Item(164,19).Value2 = "09/07/2023 23:11"
$ws.Cells.Item(164,20).Value2 = 3.07
$ws.Cells.Item(164,21).Value2 = "15/07/2023 02:29"
$ws.Cells.Item(164,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/atletico-go-sport-recife/KlNQMIY6/"

# rotate rows (166, 167)
$ws.Cells.Item(166,6).Value2 = "Ponte Preta"
$ws.Cells.Item(166,7).Value2 = 0
$ws.Cells.Item(166,8).Value2 = "Tombense"
$ws.Cells.Item(166,9).Value2 = 1
$ws.Cells.Item(166,10).Value2 = 2.01
$ws.Cells.Item(166,11).Value2 = "08/07/2023 23:12"
$ws.Cells.Item(166,12).Value2 = 2.2
$ws.Cells.Item(166,13).Value2 = "15/07/2023 21:57"
$ws.Cells.Item(166,14).Value2 = 3.2
$ws.Cells.Item(166,15).Value2 = "08/07/2023 23:12"
$ws.Cells.Item(166,16).Value2 = 3.14
$ws.Cells.Item(166,17).Value2 = "15/07/2023 21:57"
$ws.Cells.Item(166,18).Value2 = 4.29
$ws.Cells.Item(166,19).Value2 = "08/07/2023 23:12"
$ws.Cells.Item(166,20).Value2 = 3.83
$ws.Cells.Item(166,21).Value2 = "15/07/2023 21:57"
$ws.Cells.Item(166,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/ponte-preta-tombense/xQXLNbJ0/"
$ws.Cells.Item(167,6).Value2 = "Londrina"
$ws.Cells.Item(167,7).Value2 = 1
$ws.Cells.Item(167,8).Value2 = "Vila Nova FC"
$ws.Cells.Item(167,9).Value2 = 0
$ws.Cells.Item(167,10).Value2 = 3.53
$ws.Cells.Item(167,11).Value2 = "11/07/2023 01:11"
$ws.Cells.Item(167,12).Value2 = 4.33
$ws.Cells.Item(167,13).Value2 = "15/07/2023 21:56"
$ws.Cells.Item(167,14).Value2 = 3.06
$ws.Cells.Item(167,15).Value2 = "11/07/2023 01:11"
$ws.Cells.Item(167,16).Value2 = 3.17
$ws.Cells.Item(167,17).Value2 = "15/07/2023 21:56"
$ws.Cells.Item(167,18).Value2 = 2.21
$ws.Cells.Item(167,19).Value2 = "11/07/2023 01:11"
$ws.Cells.Item(167,20).Value2 = 2.05
$ws.Cells.Item(167,21).Value2 = "15/07/2023 21:56"
$ws.Cells.Item(167,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/londrina-vila-nova-fc/QqZdPBf5/"

# rotate rows (177, 178, 179)
$ws.Cells.Item(177,6).Value2 = "Avai"
$ws.Cells.Item(177,7).Value2 = 2
$ws.Cells.Item(177,8).Value2 = "Sampaio Correa"
$ws.Cells.Item(177,9).Value2 = 0
$ws.Cells.Item(177,10).Value2 = 2.2
$ws.Cells.Item(177,11).Value2 = "15/07/2023 02:42"
$ws.Cells.Item(177,12).Value2 = 2.05
$ws.Cells.Item(177,13).Value2 = "20/07/2023 02:29"
$ws.Cells.Item(177,14).Value2 = 3.05
$ws.Cells.Item(177,15).Value2 = "15/07/2023 02:42"
$ws.Cells.Item(177,16).Value2 = 2.98
$ws.Cells.Item(177,17).Value2 = "20/07/2023 02:29"
$ws.Cells.Item(177,18).Value2 = 3.82
$ws.Cells.Item(177,19).Value2 = "15/07/2023 02:42"
$ws.Cells.Item(177,20).Value2 = 4.71
$ws.Cells.Item(177,21).Value2 = "20/07/2023 02:29"
$ws.Cells.Item(177,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/avai-sampaio-correa/fPQDLTgU/"
$ws.Cells.Item(178,6).Value2 = "ABC"
$ws.Cells.Item(178,7).Value2 = 0
$ws.Cells.Item(178,8).Value2 = "Guarani"
$ws.Cells.Item(178,9).Value2 = 1
$ws.Cells.Item(178,10).Value2 = 2.6
$ws.Cells.Item(178,11).Value2 = "15/07/2023 16:12"
$ws.Cells.Item(178,12).Value2 = 2.6
$ws.Cells.Item(178,13).Value2 = "20/07/2023 02:29"
$ws.Cells.Item(178,14).Value2 = 2.85
$ws.Cells.Item(178,15).Value2 = "15/07/2023 16:12"
$ws.Cells.Item(178,16).Value2 = 2.81
$ws.Cells.Item(178,17).Value2 = "20/07/2023 02:26"
$ws.Cells.Item(178,18).Value2 = 3.08
$ws.Cells.Item(178,19).Value2 = "15/07/2023 16:12"
$ws.Cells.Item(178,20).Value2 = 3.42
$ws.Cells.Item(178,21).Value2 = "20/07/2023 02:29"
$ws.Cells.Item(178,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/abc-guarani/KS9yy9gH/"
$ws.Cells.Item(179,6).Value2 = "Sport Recife"
$ws.Cells.Item(179,7).Value2 = 1
$ws.Cells.Item(179,8).Value2 = "Vitoria"
$ws.Cells.Item(179,9).Value2 = 2
$ws.Cells.Item(179,10).Value2 = 1.69
$ws.Cells.Item(179,11).Value2 = "16/07/2023 23:12"
$ws.Cells.Item(179,12).Value2 = 1.74
$ws.Cells.Item(179,13).Value2 = "20/07/2023 01:57"
$ws.Cells.Item(179,14).Value2 = 3.51
$ws.Cells.Item(179,15).Value2 = "16/07/2023 23:12"
$ws.Cells.Item(179,16).Value2 = 3.42
$ws.Cells.Item(179,17).Value2 = "20/07/2023 01:54"
$ws.Cells.Item(179,18).Value2 = 5.8
$ws.Cells.Item(179,19).Value2 = "16/07/2023 23:12"
$ws.Cells.Item(179,20).Value2 = 5.96
$ws.Cells.Item(179,21).Value2 = "20/07/2023 01:57"
$ws.Cells.Item(179,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/sport-recife-vitoria/SCeSvifh/"

# rotate rows (203, 204)
$ws.Cells.Item(203,6).Value2 = "Criciuma"
$ws.Cells.Item(203,7).Value2 = 2
$ws.Cells.Item(203,8).Value2 = "Ponte Preta"
$ws.Cells.Item(203,9).Value2 = 1
$ws.Cells.Item(203,10).Value2 = 1.67
$ws.Cells.Item(203,11).Value2 = "30/07/2023 23:12"
$ws.Cells.Item(203,12).Value2 = 1.73
$ws.Cells.Item(203,13).Value2 = "02/08/2023 23:59"
$ws.Cells.Item(203,14).Value2 = 3.34
$ws.Cells.Item(203,15).Value2 = "30/07/2023 23:12"
$ws.Cells.Item(203,16).Value2 = 3.42
$ws.Cells.Item(203,17).Value2 = "02/08/2023 23:59"
$ws.Cells.Item(203,18).Value2 = 5.96
$ws.Cells.Item(203,19).Value2 = "30/07/2023 23:12"
$ws.Cells.Item(203,20).Value2 = 6.01
$ws.Cells.Item(203,21).Value2 = "02/08/2023 23:59"
$ws.Cells.Item(203,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/criciuma-ponte-preta/tKsyfNad/"
$ws.Cells.Item(204,6).Value2 = "Ituano"
$ws.Cells.Item(204,7).Value2 = 1
$ws.Cells.Item(204,8).Value2 = "Tombense"
$ws.Cells.Item(204,9).Value2 = 0
$ws.Cells.Item(204,10).Value2 = 2.11
$ws.Cells.Item(204,11).Value2 = "30/07/2023 16:12"
$ws.Cells.Item(204,12).Value2 = 2.14
$ws.Cells.Item(204,13).Value2 = "02/08/2023 23:58"
$ws.Cells.Item(204,14).Value2 = 3.17
$ws.Cells.Item(204,15).Value2 = "30/07/2023 16:12"
$ws.Cells.Item(204,16).Value2 = 2.98
$ws.Cells.Item(204,17).Value2 = "02/08/2023 23:58"
$ws.Cells.Item(204,18).Value2 = 3.68
$ws.Cells.Item(204,19).Value2 = "30/07/2023 16:12"
$ws.Cells.Item(204,20).Value2 = 4.31
$ws.Cells.Item(204,21).Value2 = "02/08/2023 23:58"
$ws.Cells.Item(204,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/ituano-tombense/fNVII6xc/"

# rotate rows (205, 206)
$ws.Cells.Item(205,6).Value2 = "Juventude"
$ws.Cells.Item(205,7).Value2 = 1
$ws.Cells.Item(205,8).Value2 = "Novorizontino"
$ws.Cells.Item(205,9).Value2 = 0
$ws.Cells.Item(205,10).Value2 = 2.38
$ws.Cells.Item(205,11).Value2 = "30/07/2023 20:42"
$ws.Cells.Item(205,12).Value2 = 2.46
$ws.Cells.Item(205,13).Value2 = "03/08/2023 02:21"
$ws.Cells.Item(205,14).Value2 = 2.96
$ws.Cells.Item(205,15).Value2 = "30/07/2023 20:42"
$ws.Cells.Item(205,16).Value2 = 2.94
$ws.Cells.Item(205,17).Value2 = "03/08/2023 02:28"
$ws.Cells.Item(205,18).Value2 = 3.51
$ws.Cells.Item(205,19).Value2 = "30/07/2023 20:42"
$ws.Cells.Item(205,20).Value2 = 3.48
$ws.Cells.Item(205,21).Value2 = "03/08/2023 02:21"
$ws.Cells.Item(205,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/esporte-clube-juventude-novorizontino/fLjabxir/"
$ws.Cells.Item(206,6).Value2 = "Vila Nova FC"
$ws.Cells.Item(206,7).Value2 = 0
$ws.Cells.Item(206,8).Value2 = "Sport Recife"
$ws.Cells.Item(206,9).Value2 = 1
$ws.Cells.Item(206,10).Value2 = 2.3
$ws.Cells.Item(206,11).Value2 = "30/07/2023 20:42"
$ws.Cells.Item(206,12).Value2 = 2.31
$ws.Cells.Item(206,13).Value2 = "03/08/2023 02:19"
$ws.Cells.Item(206,14).Value2 = 2.97
$ws.Cells.Item(206,15).Value2 = "30/07/2023 20:42"
$ws.Cells.Item(206,16).Value2 = 2.9
$ws.Cells.Item(206,17).Value2 = "03/08/2023 02:19"
$ws.Cells.Item(206,18).Value2 = 3.68
$ws.Cells.Item(206,19).Value2 = "30/07/2023 20:42"
$ws.Cells.Item(206,20).Value2 = 3.9
$ws.Cells.Item(206,21).Value2 = "03/08/2023 02:27"
$ws.Cells.Item(206,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/vila-nova-fc-sport-recife/2TrXf3pj/"

# append new row 282 (same formatting as the previous last data row, 281)
$ws.Cells.Item(281,1).Copy()
$ws.Cells.Item(282,1).PasteSpecial(-4122)
$ws.Cells.Item(281,5).Copy()
$ws.Cells.Item(282,5).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(282,1).Value2 = 281
$ws.Cells.Item(282,2).Value2 = "brazil"
$ws.Cells.Item(282,3).Value2 = "serie-b"
$ws.Cells.Item(1000,26).Formula = "=TEXT(2023,""0"")"
$ws.Cells.Item(1000,26).Copy()
$ws.Cells.Item(282,4).PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Cells.Item(1000,26).ClearContents()
$ws.Cells.Item(282,5).Value2 = 45189.10416666666
$ws.Cells.Item(282,6).Value2 = "Sampaio Correa"
$ws.Cells.Item(282,7).Value2 = 2
$ws.Cells.Item(282,8).Value2 = "Vila Nova FC"
$ws.Cells.Item(282,9).Value2 = 1
$ws.Cells.Item(282,10).Value2 = 2.91
$ws.Cells.Item(282,11).Value2 = "16/09/2023 21:12"
$ws.Cells.Item(282,12).Value2 = 2.98
$ws.Cells.Item(282,13).Value2 = "20/09/2023 02:23"
$ws.Cells.Item(282,14).Value2 = 2.75
$ws.Cells.Item(282,15).Value2 = "16/09/2023 21:12"
$ws.Cells.Item(282,16).Value2 = 2.75
$ws.Cells.Item(282,17).Value2 = "20/09/2023 02:08"
$ws.Cells.Item(282,18).Value2 = 2.98
$ws.Cells.Item(282,19).Value2 = "16/09/2023 21:12"
$ws.Cells.Item(282,20).Value2 = 3
$ws.Cells.Item(282,21).Value2 = "20/09/2023 02:23"
$ws.Cells.Item(282,22).Value2 = "https://www.betexplorer.com/football/brazil/serie-b/sampaio-correa-vila-nova-fc/UXNX5rBs/"

$wb.Save()
